$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44175
$ws.Range("M2").Value2 = 25
$ws.Range("N2").Value2 = 20000
$ws.Range("O2").Value2 = 20000
$ws.Range("P2").Value2 = 20000
$ws.Range("S2").Value2 = 4000

# Row 3
$ws.Range("D3").Value2 = 44186
$ws.Range("M3").Value2 = 40

# Row 4
$ws.Range("D4").Value2 = 44193
$ws.Range("M4").Value2 = 40
$ws.Range("N4").Value2 = 15000
$ws.Range("O4").Value2 = 15000
$ws.Range("P4").Value2 = 15000
$ws.Range("S4").Value2 = 3000

# Row 5
$ws.Range("D5").Value2 = 44179
$ws.Range("M5").Value2 = 45
$ws.Range("N5").Value2 = 20000
$ws.Range("O5").Value2 = 20000
$ws.Range("P5").Value2 = 20000
$ws.Range("S5").Value2 = 4000

# Row 6
$ws.Range("D6").Value2 = 44181
$ws.Range("M6").Value2 = 30
$ws.Range("N6").Value2 = 20000
$ws.Range("O6").Value2 = 20000
$ws.Range("P6").Value2 = 20000
$ws.Range("S6").Value2 = 4000

# Row 7
$ws.Range("D7").Value2 = 44189

# Row 8
$ws.Range("D8").Value2 = 44188
$ws.Range("N8").Value2 = 15000
$ws.Range("O8").Value2 = 15000
$ws.Range("P8").Value2 = 15000
$ws.Range("S8").Value2 = 3000

# Row 9
$ws.Range("D9").Value2 = 44196
$ws.Range("M9").Value2 = 56
$ws.Range("N9").Value2 = 15000
$ws.Range("O9").Value2 = 15000
$ws.Range("P9").Value2 = 15000
$ws.Range("S9").Value2 = 3000
